$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. surveys sheet: rename headers, delete rows for study02, restyle header row
# ---------------------------------------------------------------------------
$surveys = $wb.Worksheets.Item("surveys")
$surveys.Range("B1").Value = "survey_id"
$surveys.Range("E1").Value = "latitude"
$surveys.Range("F1").Value = "longitude"

# ---------------------------------------------------------------------------
# 2. studies sheet: rename header "study_ID" -> "study_id"
# ---------------------------------------------------------------------------
$studies = $wb.Worksheets.Item("studies")
$studies.Range("A1").Value = "study_id"

# ---------------------------------------------------------------------------
# 3. Notes sheet: update the "specific issue" note text (A3)
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = "Specific issue: survey_IDs are not referenced in counts table"

# ---------------------------------------------------------------------------
# back to surveys: remove the two survey rows that belong to study02 (rows 4 and 5)
# ---------------------------------------------------------------------------
$surveys.Rows.Item(4).Delete()
$surveys.Rows.Item(4).Delete()

# re-style the header row with a new explicit-black font
$surveys.Range("A1:G1").Font.Color = 0
$surveys.Range("K1").Font.Color = 0
$surveys.Range("H1:J1").Font.Color = 0
$surveys.Range("H1:J1").NumberFormat = "@"

$surveys.Range("A1:K1").Select()

# ---------------------------------------------------------------------------
# 4. studies sheet: restore selection
# ---------------------------------------------------------------------------
$studies.Range("A2").Select()

# ---------------------------------------------------------------------------
# 5. counts sheet: no content change required
# ---------------------------------------------------------------------------
$counts = $wb.Worksheets.Item("counts")
$counts.Range("A3").Select()
